$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 283 (existing rows 283:302 shift down to 285:304)
$ws.Rows.Item(283).Insert()
$ws.Rows.Item(283).Insert()

# ---- Populate new row 283 ----
$ws.Range("A283").Value = 11
$ws.Range("B283").Value = "Vega Monumental Concepción"
$ws.Range("C283").Value = "Bíobío"
$ws.Range("D283").Value = 45106
$ws.Range("D283").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E283").Value = 8
$ws.Range("F283").Value = "Fruta"
$ws.Range("G283").Value = 100101
$ws.Range("H283").Value = "Berries"
$ws.Range("I283").Value = 100101007
$ws.Range("J283").Value = "Kiwi"
$ws.Range("K283").Value = "Hayward"
$ws.Range("L283").Value = "Primera"
$ws.Range("M283").Value = 120
$ws.Range("N283").Value = 10000
$ws.Range("O283").Value = 10000
$ws.Range("P283").Value = 10000
$ws.Range("Q283").Value = "$/bandeja 18 kilos"
$ws.Range("R283").Value = "Región de O'Higgins"
$ws.Range("S283").Value = 556
$ws.Range("T283").Value = 18

# ---- Populate new row 284 ----
$ws.Range("A284").Value = 11
$ws.Range("B284").Value = "Vega Monumental Concepción"
$ws.Range("C284").Value = "Bíobío"
$ws.Range("D284").Value = 45106
$ws.Range("D284").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E284").Value = 8
$ws.Range("F284").Value = "Fruta"
$ws.Range("G284").Value = 100101
$ws.Range("H284").Value = "Berries"
$ws.Range("I284").Value = 100101007
$ws.Range("J284").Value = "Kiwi"
$ws.Range("K284").Value = "Hayward"
$ws.Range("L284").Value = "Segunda"
$ws.Range("M284").Value = 100
$ws.Range("N284").Value = 8000
$ws.Range("O284").Value = 8000
$ws.Range("P284").Value = 8000
$ws.Range("Q284").Value = "$/bandeja 18 kilos"
$ws.Range("R284").Value = "Región de O'Higgins"
$ws.Range("S284").Value = 444
$ws.Range("T284").Value = 18
